$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing color for row 6 ---
$ws.Range("K6").Value = "Green"

# --- Row 9 ---
$ws.Range("A9").Value = 28252
$ws.Range("B9").Value = "American Gouldian"
$ws.Range("C9").Value = "North America"
$ws.Range("D9").Value = "15A"
$ws.Range("E9").Value = "Male"
$ws.Range("F9").Value = 223
$ws.Range("G9").Value = 111
$ws.Range("H9").Value = "15/05/2023"
$ws.Range("I9").Value = "Red"
$ws.Range("J9").Value = "Purple"
$ws.Range("K9").Value = "Pastel"

# --- Row 10 ---
$ws.Range("A10").Value = 2222
$ws.Range("B10").Value = "European Gouldian"
$ws.Range("C10").Value = "Western Europe"
$ws.Range("D10").Value = 1111
$ws.Range("E10").Value = "Male"
$ws.Range("F10").Value = 223
$ws.Range("G10").Value = 111
$ws.Range("H10").Value = "15/05/2023"
$ws.Range("I10").Value = "Red"
$ws.Range("J10").Value = "Purple"
$ws.Range("K10").Value = "Pastel"

# --- Row 11 ---
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "American Gouldian"
$ws.Range("C11").Value = "Central America"
$ws.Range("D11").Value = "12A"
$ws.Range("E11").Value = "Female"
$ws.Range("F11").Value = 223
$ws.Range("G11").Value = 111
$ws.Range("H11").Value = "15/05/2023"
$ws.Range("I11").Value = "Black"
$ws.Range("J11").Value = "Purple"
$ws.Range("K11").Value = "Green"

# --- Row 12 ---
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "American Gouldian"
$ws.Range("C12").Value = "North America"
$ws.Range("D12").Value = "12A"
$ws.Range("E12").Value = "Female"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 111
$ws.Range("H12").Value = "15/05/2023"
$ws.Range("I12").Value = "Black"
$ws.Range("J12").Value = "Purple"
$ws.Range("K12").Value = "Green"

# --- Row 13 ---
$ws.Range("A13").Value = 555
$ws.Range("B13").Value = "American Gouldian"
$ws.Range("C13").Value = "North America"
$ws.Range("D13").Value = 11
$ws.Range("E13").Value = "Female"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 435345
$ws.Range("H13").Value = "15/05/2023"
$ws.Range("I13").Value = "Red"
$ws.Range("J13").Value = "Purple"

# --- Update selection to match the final workbook state ---
$ws.Range("L6").Select()
